# Auto-generated: applies scheduled market-data refresh to Garuda_Profits sheets
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 91001730
$ws.Range("I86").Value = 112944.78
$ws.Range("J86").Value = 500001250
$ws.Range("K86").Value = 112944.78
$ws.Range("L86").Value = 500001250
$ws.Range("M86").Value = -111821.78
$ws.Range("N86").Value = -500003496
$ws.Range("H89").Value = 91001730
$ws.Range("I89").Value = 112944.78
$ws.Range("J89").Value = 500001250
$ws.Range("K89").Value = 564723.9
$ws.Range("L89").Value = 2500006250
$ws.Range("M89").Value = -559107.9
$ws.Range("N89").Value = -2500017482
$ws.Range("H111").Value = 2560.9
$ws.Range("I111").Value = 2879.8333
$ws.Range("J111").Value = 2082.5
$ws.Range("K111").Value = 8639.499899999999
$ws.Range("L111").Value = 6247.5
$ws.Range("M111").Value = -5572.499899999999
$ws.Range("N111").Value = -12381.5
$ws.Range("H125").Value = 250001860
$ws.Range("J125").Value = 2466.6667
$ws.Range("L125").Value = 22200.0003
$ws.Range("N125").Value = -27120.0003
$ws.Range("H137").Value = 18519756
$ws.Range("I137").Value = 1087.0513
$ws.Range("J137").Value = 66668292
$ws.Range("K137").Value = 3261.1539
$ws.Range("L137").Value = 200004876
$ws.Range("M137").Value = -711.1539000000002
$ws.Range("N137").Value = -200009976

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28116.258
$ws.Range("I32").Value = 30283.178
$ws.Range("J32").Value = 20615.385
$ws.Range("K32").Value = 30283.178
$ws.Range("L32").Value = 20615.385
$ws.Range("M32").Value = -29996.178
$ws.Range("N32").Value = -21189.385
$ws.Range("H61").Value = 1843.3214
$ws.Range("I61").Value = 1460.826
$ws.Range("J61").Value = 3602.8
$ws.Range("K61").Value = 1460.826
$ws.Range("L61").Value = 3602.8
$ws.Range("M61").Value = -1248.826
$ws.Range("N61").Value = -4026.8
$ws.Range("H74").Value = 958.1087
$ws.Range("I74").Value = 892.2143
$ws.Range("J74").Value = 1650
$ws.Range("K74").Value = 892.2143
$ws.Range("L74").Value = 1650
$ws.Range("M74").Value = -18.21429999999998
$ws.Range("N74").Value = -3398
$ws.Range("H77").Value = 958.1087
$ws.Range("I77").Value = 892.2143
$ws.Range("J77").Value = 1650
$ws.Range("K77").Value = 4461.0715
$ws.Range("L77").Value = 8250
$ws.Range("M77").Value = -93.07150000000001
$ws.Range("N77").Value = -16986
$ws.Range("H136").Value = 1843.3214
$ws.Range("I136").Value = 1460.826
$ws.Range("J136").Value = 3602.8
$ws.Range("K136").Value = 4382.478
$ws.Range("L136").Value = 10808.4
$ws.Range("M136").Value = -1832.478
$ws.Range("N136").Value = -15908.4

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 378.33334
$ws.Range("I80").Value = 429.5
$ws.Range("J80").Value = 359.72726
$ws.Range("K80").Value = 429.5
$ws.Range("L80").Value = 359.72726
$ws.Range("M80").Value = 568.5
$ws.Range("N80").Value = -2355.72726
$ws.Range("H83").Value = 378.33334
$ws.Range("I83").Value = 429.5
$ws.Range("J83").Value = 359.72726
$ws.Range("K83").Value = 2147.5
$ws.Range("L83").Value = 1798.6363
$ws.Range("M83").Value = 2844.5
$ws.Range("N83").Value = -11782.6363
$ws.Range("H86").Value = 2926.1667
$ws.Range("I86").Value = 3350
$ws.Range("J86").Value = 2502.3333
$ws.Range("K86").Value = 3350
$ws.Range("L86").Value = 2502.3333
$ws.Range("M86").Value = -2227
$ws.Range("N86").Value = -4748.3333
$ws.Range("H89").Value = 2926.1667
$ws.Range("I89").Value = 3350
$ws.Range("J89").Value = 2502.3333
$ws.Range("K89").Value = 16750
$ws.Range("L89").Value = 12511.6665
$ws.Range("M89").Value = -11134
$ws.Range("N89").Value = -23743.6665
$ws.Range("H94").Value = 374.75
$ws.Range("I94").Value = 359.81818
$ws.Range("K94").Value = 359.81818
$ws.Range("M94").Value = 91.18182000000002
$ws.Range("H99").Value = 747.7778
$ws.Range("I99").Value = 718.5714
$ws.Range("J99").Value = 850
$ws.Range("K99").Value = 718.5714
$ws.Range("L99").Value = 850
$ws.Range("M99").Value = 779.4286
$ws.Range("N99").Value = -3846
$ws.Range("H105").Value = 1965.7142
$ws.Range("I105").Value = 1830
$ws.Range("J105").Value = 2101.4285
$ws.Range("K105").Value = 1830
$ws.Range("L105").Value = 2101.4285
$ws.Range("M105").Value = -83
$ws.Range("N105").Value = -5595.4285
$ws.Range("H134").Value = 5152.077
$ws.Range("I134").Value = 5909.7144
$ws.Range("J134").Value = 3223.5454
$ws.Range("K134").Value = 17729.1432
$ws.Range("L134").Value = 9670.636200000001
$ws.Range("M134").Value = -15194.1432
$ws.Range("N134").Value = -14740.6362

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4368.1724
$ws.Range("I31").Value = 2964.611
$ws.Range("J31").Value = 6664.909
$ws.Range("K31").Value = 2964.611
$ws.Range("L31").Value = 6664.909
$ws.Range("M31").Value = -2669.611
$ws.Range("N31").Value = -7254.909
$ws.Range("H34").Value = 4368.1724
$ws.Range("I34").Value = 2964.611
$ws.Range("J34").Value = 6664.909
$ws.Range("K34").Value = 2964.611
$ws.Range("L34").Value = 6664.909
$ws.Range("M34").Value = -2762.611
$ws.Range("N34").Value = -7068.909
$ws.Range("H105").Value = 729.8333
$ws.Range("I105").Value = 747.5
$ws.Range("K105").Value = 747.5
$ws.Range("M105").Value = 999.5
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 692.8570999999999
$ws.Range("I4").Value = 170
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 510
$ws.Range("L4").Value = 6000
$ws.Range("M4").Value = -398
$ws.Range("N4").Value = -6224
$ws.Range("H107").Value = 394.3889
$ws.Range("I107").Value = 272.36365
$ws.Range("J107").Value = 586.1429000000001
$ws.Range("K107").Value = 817.09095
$ws.Range("L107").Value = 1758.4287
$ws.Range("M107").Value = 1102.90905
$ws.Range("N107").Value = -5598.4287

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2940.5557
$ws.Range("I40").Value = 2923.5715
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2923.5715
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -2787.5715
$ws.Range("N40").Value = -3272
$ws.Range("H61").Value = 11905903
$ws.Range("I61").Value = 1002.1429
$ws.Range("J61").Value = 23810804
$ws.Range("K61").Value = 1002.1429
$ws.Range("L61").Value = 23810804
$ws.Range("M61").Value = -800.1429000000001
$ws.Range("N61").Value = -23811208
$ws.Range("H82").Value = 1741.7142
$ws.Range("I82").Value = 1741.7142
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1741.7142
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1380.7142
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 1741.7142
$ws.Range("I85").Value = 1741.7142
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1741.7142
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -493.7141999999999
$ws.Range("N85").ClearContents()
$ws.Range("H113").Value = 11905903
$ws.Range("I113").Value = 1002.1429
$ws.Range("J113").Value = 23810804
$ws.Range("K113").Value = 1002.1429
$ws.Range("L113").Value = 23810804
$ws.Range("M113").Value = 1167.8571
$ws.Range("N113").Value = -23815144

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3924.25
$ws.Range("I81").Value = 3924.25
$ws.Range("K81").Value = 7848.5
$ws.Range("M81").Value = -6787.5
$ws.Range("H84").Value = 3924.25
$ws.Range("I84").Value = 3924.25
$ws.Range("K84").Value = 39242.5
$ws.Range("M84").Value = -33938.5
$ws.Range("H126").Value = 33339536
$ws.Range("J126").Value = 4365.2856
$ws.Range("L126").Value = 13095.8568
$ws.Range("N126").Value = -18035.8568
